# Fruta / hortaliza, semanal
# A new weekly observation is inserted as row 57, pushing the existing
# rows 57-89 down to 58-90 (dimension grows from A1:T89 to A1:T90).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 57; everything below shifts down one row.
$ws.Rows("57:57").Insert()

# Populate the newly inserted row with the new observation's data.
$ws.Cells.Item(57, 1).Value  = 4
$ws.Cells.Item(57, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(57, 3).Value  = "Los Lagos"
$ws.Cells.Item(57, 4).Value  = 44489
$ws.Cells.Item(57, 5).Value  = 10
$ws.Cells.Item(57, 6).Value  = "Fruta"
$ws.Cells.Item(57, 7).Value  = 100108
$ws.Cells.Item(57, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(57, 9).Value  = 100108002
$ws.Cells.Item(57, 10).Value = "Mango"
$ws.Cells.Item(57, 11).Value = "Sin especificar"
$ws.Cells.Item(57, 12).Value = "Primera"
$ws.Cells.Item(57, 13).Value = 60
$ws.Cells.Item(57, 14).Value = 7500
$ws.Cells.Item(57, 15).Value = 8000
$ws.Cells.Item(57, 16).Value = 7750
$ws.Cells.Item(57, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(57, 18).Value = "Perú"
$ws.Cells.Item(57, 19).Value = 1938
$ws.Cells.Item(57, 20).Value = 4
